# Correccion en el calculo de las fechas de firma
# Updates the EIXIDA/TORNADA dates & hours for the first two trip rows,
# clears the next three trip rows entirely, and fixes the truncated
# "INSP. JEFA, JEFA PROVINCIAL ..." signature line.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Row 12 (1-based) : 01-08-2025 / 07:00 -> 01-08-2025 / 15:00  =>  13-08-2025 / 23:00 -> 14-08-2025 / 07:00
$row = $t.Rows.Item(12)
$row.Cells.Item(3).Range.Text  = "13-08-2025"
$row.Cells.Item(4).Range.Text  = "23:00"
$row.Cells.Item(5).Range.Text  = "14-08-2025"
$row.Cells.Item(6).Range.Text  = "07:00"

# --- Row 13 (1-based) : 02-08-2025 / 07:00 -> 02-08-2025 / 15:00  =>  14-08-2025 / 15:00 -> 14-08-2025 / 23:00
$row = $t.Rows.Item(13)
$row.Cells.Item(3).Range.Text  = "14-08-2025"
$row.Cells.Item(4).Range.Text  = "15:00"
$row.Cells.Item(5).Range.Text  = "14-08-2025"
$row.Cells.Item(6).Range.Text  = "23:00"

# --- Rows 14, 15, 16 (1-based) : wipe the whole trip entry (itinerary, dates, hours, V.P., km)
foreach ($r in 14,15,16) {
    $row = $t.Rows.Item($r)
    $row.Cells.Item(1).Range.Text  = ""
    $row.Cells.Item(3).Range.Text  = ""
    $row.Cells.Item(4).Range.Text  = ""
    $row.Cells.Item(5).Range.Text  = ""
    $row.Cells.Item(6).Range.Text  = ""
    $row.Cells.Item(10).Range.Text = ""
    $row.Cells.Item(11).Range.Text = ""
}

# --- Fix the truncated approval line
foreach ($p in $d.Content.Paragraphs) {
    if ($p.Range.Text -like "*INSP. JEFA, JEFA PROVINCIAL UN*") {
        $p.Range.Text = "INSP. JEFA, JEFA PROVINCIAL Unidad Adscrita a la P.N."
        break
    }
}
